# EstrategiaEntrenamiento.xlsx edit script
# Implements: new "Opcion 2" / Validation+Holdout block gets a "Final Train" row
# inserted on Hoja2, and the workbook window / selection metadata is refreshed.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Hoja1")
$ws2 = $wb.Worksheets.Item("Hoja2")

# ---------------------------------------------------------------------------
# Hoja2: insert 3 new rows right below the "Opcion 2" / Train block (old row 9)
# so the merged A7:A9 label area grows to A7:A12 and everything further down
# the sheet shifts down by 3 rows (old row 11 -> 14, old row 14 -> 17, ...).
# ---------------------------------------------------------------------------
$ws2.Rows("9:11").Insert()

# Move the "Validation" caption up into row 8 (it used to sit in row 9,
# which is now row 12 after the insert above).
$ws2.Range("N12").Cut($ws2.Range("N8"))
$ws2.Range("N12").ClearFormats()

# New "Final Train" row: C10:N10 share the same new fill color (a fresh
# olive/yellow, not previously in the palette) but are NOT merged - each
# cell keeps the style individually, matching the other label cells.
$ws2.Range("C10:N10").Interior.Color = 52428
$ws2.Range("C10").Value = "Final Train"

# Move "Holdout" down into row 11 (it used to sit beside Validation in the
# old row 9, now pushed to row 12 by the insert).
$ws2.Range("P12").Cut($ws2.Range("P11"))
$ws2.Range("P12").ClearFormats()

# Refresh the sheet's selection / scroll anchor to the new "Final Train" row.
$ws2.Range("N10").Select()

# ---------------------------------------------------------------------------
# Workbook-level cosmetic metadata refresh (mirrors a normal save-after-edit
# from a new working copy of the repo).
# ---------------------------------------------------------------------------
$excel.Windows.Item(1).Left = -120

$wb.Save()
